# Generate Report for Handoff
#
# Adds two new handoff entries ("2d6d42a7-817a-4ea7-89a8-ab7fd7857102" and
# "66f8ed3e-b652-4472-84d9-8a2a16539b78") around the existing
# "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c" entry (which moves from row 5 to
# row 6) on all three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Step 1: insert two new (style-preserving) rows before the old row 5 on
# each sheet, by copying row 5 (which already carries the correct cell
# styles for the new rows) and inserting the copy twice. This pushes the
# pre-existing "37ea1dd4" row from row 5 down to row 7, after which we
# will overwrite the cell values of rows 5, 6 and 7 with the final data.
# ---------------------------------------------------------------------
foreach ($ws in @($wsOverview, $wsZhCn, $wsDeDe)) {
    $ws.Rows.Item(5).Copy()
    $ws.Rows.Item(5).Insert()
    $ws.Rows.Item(5).Copy()
    $ws.Rows.Item(5).Insert()
}

# ---------------------------------------------------------------------
# Step 2: overwrite cell values for rows 5-7 (Overview sheet)
#   row5 = 2d6d42a7-817a-4ea7-89a8-ab7fd7857102 (new)
#   row6 = 37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c (previously row5, same data)
#   row7 = 66f8ed3e-b652-4472-84d9-8a2a16539b78 (new)
# ---------------------------------------------------------------------
$wsOverview.Range("A5").Value = "2d6d42a7-817a-4ea7-89a8-ab7fd7857102.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-26-11 18:26:49"

$wsOverview.Range("A6").Value = "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md"
$wsOverview.Range("B6").Value = "Ready for handoff"
$wsOverview.Range("C6").Value = "Ready for handoff"
$wsOverview.Range("D6").Value = "2016-25-11 18:25:12"

$wsOverview.Range("A7").Value = "66f8ed3e-b652-4472-84d9-8a2a16539b78.md"
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-26-11 18:26:49"

# ---------------------------------------------------------------------
# Step 3: overwrite cell values for rows 5-7 (zh-cn sheet)
# ---------------------------------------------------------------------
$wsZhCn.Range("A5").Value = "2d6d42a7-817a-4ea7-89a8-ab7fd7857102.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "2d6d42a7-817a-4ea7-89a8-ab7fd7857102.79bf10be3a863fe351819c8b8fe60ba7cad90ac2.zh-cn.xlf"
$wsZhCn.Range("E5").Value = "2016-03-11 18:26:46"
$wsZhCn.Range("H5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I5").Value = "Include"

$wsZhCn.Range("A6").Value = "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md"
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.zh-cn.xlf"
$wsZhCn.Range("E6").Value = "2016-03-11 18:25:09"
$wsZhCn.Range("H6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I6").Value = "Include"

$wsZhCn.Range("A7").Value = "66f8ed3e-b652-4472-84d9-8a2a16539b78.md"
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "66f8ed3e-b652-4472-84d9-8a2a16539b78.264f874bbf759ffd3062937f08dc13c31250895d.zh-cn.xlf"
$wsZhCn.Range("E7").Value = "2016-03-11 18:26:46"
$wsZhCn.Range("H7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I7").Value = "Include"

# ---------------------------------------------------------------------
# Step 4: overwrite cell values for rows 5-7 (de-de sheet)
# ---------------------------------------------------------------------
$wsDeDe.Range("A5").Value = "2d6d42a7-817a-4ea7-89a8-ab7fd7857102.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "2d6d42a7-817a-4ea7-89a8-ab7fd7857102.79bf10be3a863fe351819c8b8fe60ba7cad90ac2.de-de.xlf"
$wsDeDe.Range("E5").Value = "2016-03-11 18:26:49"
$wsDeDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I5").Value = "Include"

$wsDeDe.Range("A6").Value = "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md"
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.de-de.xlf"
$wsDeDe.Range("E6").Value = "2016-03-11 18:25:12"
$wsDeDe.Range("H6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I6").Value = "Include"

$wsDeDe.Range("A7").Value = "66f8ed3e-b652-4472-84d9-8a2a16539b78.md"
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "66f8ed3e-b652-4472-84d9-8a2a16539b78.264f874bbf759ffd3062937f08dc13c31250895d.de-de.xlf"
$wsDeDe.Range("E7").Value = "2016-03-11 18:26:49"
$wsDeDe.Range("H7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I7").Value = "Include"

# ---------------------------------------------------------------------
# Step 5: rebuild hyperlinks from scratch on every sheet, in the final
# left-to-right, top-to-bottom order, since row insertion does not shift
# pre-existing hyperlink ranges automatically in this engine.
# ---------------------------------------------------------------------
$wsOverview.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Delete()

# --- Overview sheet hyperlinks (column A only) ---
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/43f7a690a20b905f31a8c3fdf488167a3321d2e8/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/18ffabec28e7599157b81355ccfd005434c0d9e7/e2e/2a521b16-3871-45a3-90ec-45aa6e15bd71.md", "", "", "2a521b16-3871-45a3-90ec-45aa6e15bd71.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/18ffabec28e7599157b81355ccfd005434c0d9e7/e2e/4967c61c-a77c-496d-a7ea-863e9bf454f1.md", "", "", "4967c61c-a77c-496d-a7ea-863e9bf454f1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/79bf10be3a863fe351819c8b8fe60ba7cad90ac2/e2e/2d6d42a7-817a-4ea7-89a8-ab7fd7857102.md", "", "", "2d6d42a7-817a-4ea7-89a8-ab7fd7857102.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a769f066dcdfc66e2f1210d9ce9ee413c8966878/e2e/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md", "", "", "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/264f874bbf759ffd3062937f08dc13c31250895d/e2e/66f8ed3e-b652-4472-84d9-8a2a16539b78.md", "", "", "66f8ed3e-b652-4472-84d9-8a2a16539b78.md")

# --- zh-cn sheet hyperlinks (columns A, B, D, F, G) ---
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/43f7a690a20b905f31a8c3fdf488167a3321d2e8/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/43f7a690a20b905f31a8c3fdf488167a3321d2e8/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd824585b09f8bf94fa4886d3450c9c9e3636bd8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.zh-cn.xlf", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a26d046bd5f8cc9b4faea6b470cbd87ce14a33d2/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/da42f7c315701b77031b063ff74e3006fe9d43c1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.zh-cn.xlf", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.zh-cn.xlf")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/18ffabec28e7599157b81355ccfd005434c0d9e7/e2e/2a521b16-3871-45a3-90ec-45aa6e15bd71.md", "", "", "2a521b16-3871-45a3-90ec-45aa6e15bd71.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/18ffabec28e7599157b81355ccfd005434c0d9e7/e2e/2a521b16-3871-45a3-90ec-45aa6e15bd71.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f6f888c5686991894c957dcedaf56d01b1d210db/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2a521b16-3871-45a3-90ec-45aa6e15bd71.9c4f8a014caf28c40df2d99ab686e712762239f5.zh-cn.xlf", "", "", "2a521b16-3871-45a3-90ec-45aa6e15bd71.9c4f8a014caf28c40df2d99ab686e712762239f5.zh-cn.xlf")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/18ffabec28e7599157b81355ccfd005434c0d9e7/e2e/4967c61c-a77c-496d-a7ea-863e9bf454f1.md", "", "", "4967c61c-a77c-496d-a7ea-863e9bf454f1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/18ffabec28e7599157b81355ccfd005434c0d9e7/e2e/4967c61c-a77c-496d-a7ea-863e9bf454f1.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f6f888c5686991894c957dcedaf56d01b1d210db/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4967c61c-a77c-496d-a7ea-863e9bf454f1.25653456a5ae380c9fbd77dbcaa5c98ee47c440d.zh-cn.xlf", "", "", "4967c61c-a77c-496d-a7ea-863e9bf454f1.25653456a5ae380c9fbd77dbcaa5c98ee47c440d.zh-cn.xlf")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/79bf10be3a863fe351819c8b8fe60ba7cad90ac2/e2e/2d6d42a7-817a-4ea7-89a8-ab7fd7857102.md", "", "", "2d6d42a7-817a-4ea7-89a8-ab7fd7857102.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/79bf10be3a863fe351819c8b8fe60ba7cad90ac2/e2e/2d6d42a7-817a-4ea7-89a8-ab7fd7857102.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/79bf10be3a863fe351819c8b8fe60ba7cad90ac2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2d6d42a7-817a-4ea7-89a8-ab7fd7857102.79bf10be3a863fe351819c8b8fe60ba7cad90ac2.zh-cn.xlf", "", "", "2d6d42a7-817a-4ea7-89a8-ab7fd7857102.79bf10be3a863fe351819c8b8fe60ba7cad90ac2.zh-cn.xlf")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a769f066dcdfc66e2f1210d9ce9ee413c8966878/e2e/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md", "", "", "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/a769f066dcdfc66e2f1210d9ce9ee413c8966878/e2e/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/406abe3598a4c180756b3ecd0af98c86adb5b31d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.zh-cn.xlf", "", "", "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.zh-cn.xlf")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/264f874bbf759ffd3062937f08dc13c31250895d/e2e/66f8ed3e-b652-4472-84d9-8a2a16539b78.md", "", "", "66f8ed3e-b652-4472-84d9-8a2a16539b78.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/264f874bbf759ffd3062937f08dc13c31250895d/e2e/66f8ed3e-b652-4472-84d9-8a2a16539b78.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/264f874bbf759ffd3062937f08dc13c31250895d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/66f8ed3e-b652-4472-84d9-8a2a16539b78.264f874bbf759ffd3062937f08dc13c31250895d.zh-cn.xlf", "", "", "66f8ed3e-b652-4472-84d9-8a2a16539b78.264f874bbf759ffd3062937f08dc13c31250895d.zh-cn.xlf")

# --- de-de sheet hyperlinks (columns A, B, D, F, G) ---
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/43f7a690a20b905f31a8c3fdf488167a3321d2e8/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/43f7a690a20b905f31a8c3fdf488167a3321d2e8/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b1cb137ed09fbabdeb80df581d021d696143428d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.de-de.xlf", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8540ce3874af8a6ad9fbcc68008525e1f084ef6a/e2e/c003bd65-8677-4b9d-aad3-abac071d090b.md", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a0b4348416064dc28dde4c87dda48735c48d6bc7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.de-de.xlf", "", "", "c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.de-de.xlf")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/18ffabec28e7599157b81355ccfd005434c0d9e7/e2e/2a521b16-3871-45a3-90ec-45aa6e15bd71.md", "", "", "2a521b16-3871-45a3-90ec-45aa6e15bd71.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/18ffabec28e7599157b81355ccfd005434c0d9e7/e2e/2a521b16-3871-45a3-90ec-45aa6e15bd71.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/db87ba0ab29362186b7691e9956429999c4e7757/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2a521b16-3871-45a3-90ec-45aa6e15bd71.9c4f8a014caf28c40df2d99ab686e712762239f5.de-de.xlf", "", "", "2a521b16-3871-45a3-90ec-45aa6e15bd71.9c4f8a014caf28c40df2d99ab686e712762239f5.de-de.xlf")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/18ffabec28e7599157b81355ccfd005434c0d9e7/e2e/4967c61c-a77c-496d-a7ea-863e9bf454f1.md", "", "", "4967c61c-a77c-496d-a7ea-863e9bf454f1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/18ffabec28e7599157b81355ccfd005434c0d9e7/e2e/4967c61c-a77c-496d-a7ea-863e9bf454f1.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/db87ba0ab29362186b7691e9956429999c4e7757/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4967c61c-a77c-496d-a7ea-863e9bf454f1.25653456a5ae380c9fbd77dbcaa5c98ee47c440d.de-de.xlf", "", "", "4967c61c-a77c-496d-a7ea-863e9bf454f1.25653456a5ae380c9fbd77dbcaa5c98ee47c440d.de-de.xlf")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/79bf10be3a863fe351819c8b8fe60ba7cad90ac2/e2e/2d6d42a7-817a-4ea7-89a8-ab7fd7857102.md", "", "", "2d6d42a7-817a-4ea7-89a8-ab7fd7857102.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/79bf10be3a863fe351819c8b8fe60ba7cad90ac2/e2e/2d6d42a7-817a-4ea7-89a8-ab7fd7857102.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/79bf10be3a863fe351819c8b8fe60ba7cad90ac2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2d6d42a7-817a-4ea7-89a8-ab7fd7857102.79bf10be3a863fe351819c8b8fe60ba7cad90ac2.de-de.xlf", "", "", "2d6d42a7-817a-4ea7-89a8-ab7fd7857102.79bf10be3a863fe351819c8b8fe60ba7cad90ac2.de-de.xlf")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a769f066dcdfc66e2f1210d9ce9ee413c8966878/e2e/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md", "", "", "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/a769f066dcdfc66e2f1210d9ce9ee413c8966878/e2e/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff3a514f657f04f247309fad5bfb5fa8e767cbdd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.de-de.xlf", "", "", "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.de-de.xlf")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/264f874bbf759ffd3062937f08dc13c31250895d/e2e/66f8ed3e-b652-4472-84d9-8a2a16539b78.md", "", "", "66f8ed3e-b652-4472-84d9-8a2a16539b78.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/264f874bbf759ffd3062937f08dc13c31250895d/e2e/66f8ed3e-b652-4472-84d9-8a2a16539b78.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/264f874bbf759ffd3062937f08dc13c31250895d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/66f8ed3e-b652-4472-84d9-8a2a16539b78.264f874bbf759ffd3062937f08dc13c31250895d.de-de.xlf", "", "", "66f8ed3e-b652-4472-84d9-8a2a16539b78.264f874bbf759ffd3062937f08dc13c31250895d.de-de.xlf")

"Report generated for handoff: added 2d6d42a7-817a-4ea7-89a8-ab7fd7857102 and 66f8ed3e-b652-4472-84d9-8a2a16539b78" | Out-Host
